# Apply odds updates described by the diff for
# Jogos_da_Semana_FlashScore_2024-11-14.xlsx
#
# The workbook has a single worksheet with one header row (row 1) and
# match rows below it. This script updates specific odds cells in
# rows 2, 7, 8 and 10 to match the new values from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Alianza vs Bucaramanga)
$ws.Range("G2").Value = 2.5
$ws.Range("I2").Value = 3.1
$ws.Range("L2").Value = 4
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 23
$ws.Range("AK2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("AR2").Value = 81
$ws.Range("BB2").Value = 351

# Row 7 (Fortaleza vs Aguilas)
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9

# Row 8 (Ind. Medellin vs Envigado)
$ws.Range("G8").Value = 1.5
$ws.Range("Q8").Value = 1.9
$ws.Range("R8").Value = 1.95
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("W8").Value = 6.5
$ws.Range("X8").Value = 7
$ws.Range("Z8").Value = 10
$ws.Range("AE8").Value = 19
$ws.Range("AJ8").Value = 21
$ws.Range("AL8").Value = 51
$ws.Range("AN8").Value = 3.4
$ws.Range("AU8").Value = 9
$ws.Range("AY8").Value = 41
$ws.Range("BA8").Value = 151
$ws.Range("BB8").Value = 301

# Row 10 (Pereira vs La Equidad)
$ws.Range("G10").Value = 1.67
$ws.Range("I10").Value = 5.75
$ws.Range("K10").Value = 2.1
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("W10").Value = 6
$ws.Range("AC10").Value = 8
$ws.Range("AW10").Value = 7
